$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "44.592.85"
Set-TextValue "E2" "  +3.95%  "
Set-TextValue "D3" "2.433.82"
Set-TextValue "E3" "  +2.74%  "
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.09%  "
Set-TextValue "D5" "311.38"
Set-TextValue "E5" "  +3.22%  "
Set-TextValue "D6" "101.99"
Set-TextValue "E6" "  +6.21%  "
Set-TextValue "D7" "0.513"
Set-TextValue "E7" "  +1.72%  "
Set-TextValue "E8" "  -0.08%  "
Set-TextValue "E9" "  +2.84%  "
Set-TextValue "D10" "35.48"
Set-TextValue "E10" "  +3.96%  "
Set-TextValue "E11" "  +1.72%  "
Set-TextValue "E12" "  +1.10%  "
Set-TextValue "D13" "18.73"
Set-TextValue "E13" "  +2.77%  "
Set-TextValue "E14" "  +2.97%  "
Set-TextValue "D15" "2.812.43"
Set-TextValue "E15" "  +2.43%  "
Set-TextValue "D16" "2.413.47"
Set-TextValue "E16" "  +2.06%  "
Set-TextValue "E17" "  +4.65%  "
Set-TextValue "D18" "44.517.63"
Set-TextValue "E18" "  +3.81%  "
Set-TextValue "D19" "12.43"
Set-TextValue "E19" "  +2.34%  "
Set-TextValue "D20" "6.42"
Set-TextValue "E20" "  +1.90%  "
Set-TextValue "D21" "0.0₃0909"
Set-TextValue "E21" "  +2.44%  "
Set-TextValue "D22" "68.90"
Set-TextValue "E22" "  +1.31%  "
Set-TextValue "E23" "  +3.81%  "
Set-TextValue "D24" "241.27"
Set-TextValue "E24" "  +2.71%  "
Set-TextValue "D25" "2.47"
Set-TextValue "E25" "  +1.66%  "
Set-TextValue "E26" "  -0.02%  "
Set-TextValue "E27" "  +1.72%  "
Set-TextValue "E28" "  -4.23%  "
Set-TextValue "E29" "  +4.83%  "
Set-TextValue "D30" "33.41"
Set-TextValue "E30" "  +5.83%  "
Set-TextValue "D31" "48.75"
Set-TextValue "E31" "  +1.44%  "
Set-TextValue "E32" "  +16.21%  "
Set-TextValue "D33" "19.53"
Set-TextValue "E33" "  +12.14%  "
Set-TextValue "E34" "  +3.08%  "
Set-TextValue "E35" "  +0.13%  "
Set-TextValue "E36" "  +3.77%  "
Set-TextValue "E37" "  +2.82%  "
Set-TextValue "E38" "  +4.33%  "
Set-TextValue "E39" "  +4.18%  "
Set-TextValue "D40" "124.48"
Set-TextValue "E40" "  +6.89%  "
Set-TextValue "E41" "  +0.91%  "
Set-TextValue "D42" "22.01"
Set-TextValue "E42" "  -1.98%  "
Set-TextValue "D43" "2.16"
Set-TextValue "E43" "  -6.54%  "
Set-TextValue "E44" "  +3.76%  "
Set-TextValue "D45" "1.948.58"
Set-TextValue "E45" "  +0.42%  "
Set-TextValue "D46" "2.18"
Set-TextValue "E46" "  +2.34%  "
Set-TextValue "E47" "  +8.39%  "
Set-TextValue "D48" "9.81"
Set-TextValue "E48" "  +6.87%  "
Set-TextValue "D49" "1.67"
Set-TextValue "E49" "  +10.58%  "
Set-TextValue "D50" "53.47"
Set-TextValue "E50" "  +2.79%  "
Set-TextValue "D51" "73.99"
Set-TextValue "E51" "  +2.77%  "
